# Userstories 17/03 aanvullen & bug bij checkEmpySpots fixen
#
# 1. Remove the "Status" column (F) from the Userstory sheet - the
#    In Progress / Finished status tracking is no longer used.
# 2. Add three new userstory rows (1100/1200/1300) describing the
#    17/03 work: the diagonal 4-on-a-row algorithm, the AI unit-test
#    methods, and the AI & diagonal unit tests.
# 3. Move the active selection to G7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Remove the Status column entirely (was column F) ---
$ws.Columns.Item(6).Delete()

# --- Copy the existing date formatting down onto the new rows ---
$ws.Range("E11").Copy()
$ws.Range("E12:E14").PasteSpecial(-4122)

# --- New userstory rows for 17/03 ---
$ws.Range("B13").Value = "Algoritme voor diagonaal 4op een rij te vinden"
$ws.Range("C13").Value = 4
$ws.Range("D13").Value = "Lucas"
$ws.Range("E13").Value = 42080

$ws.Range("B14").Value = "Methodes voor UnitTest AI schrijven"
$ws.Range("C14").Value = 4
$ws.Range("D14").Value = "Jel"
$ws.Range("E14").Value = 42080

$ws.Range("B12").Value = "UnitTests voor AI & Diagonaal 4 op een rij te vinden"
$ws.Range("C12").Value = 4
$ws.Range("D12").Value = "Shane"
$ws.Range("E12").Value = 42080

# --- Update the selected cell on the sheet ---
$ws.Range("G7").Select()
